$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ------------------------------------------------------------------
# 1. Add the 5th column ("ACTUAL OUTCOME") and resize every column to
#    the new grid widths (this rewrites <w:tblGrid> and every <w:tcW>).
# ------------------------------------------------------------------
$t.Columns.Add() | Out-Null
$t.Columns.Item(1).Width = 70.65
$t.Columns.Item(2).Width = 99.2
$t.Columns.Item(3).Width = 85.05
$t.Columns.Item(4).Width = 106.35
$t.Columns.Item(5).Width = 92.1

# ------------------------------------------------------------------
# 2. Shade the new header cell like its neighbours.
# ------------------------------------------------------------------
$headerCell5 = $t.Cell(1, 5)
$headerCell5.Shading.Texture = 0
$headerCell5.Shading.ForegroundPatternColor = -16777216
$headerCell5.Shading.BackgroundPatternColor = 12566463

# ------------------------------------------------------------------
# 3. Fill the new header cell with "ACTUAL OUTCOME" (bold) followed by
#    the relocated _GoBack bookmark. Insert via raw XML so the bookmark
#    lands after the run instead of wrapping it.
# ------------------------------------------------------------------
$headerXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>ACTUAL OUTCOME</w:t></w:r><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>"
$headerCell5.Range.InsertXML($headerXml)

# ------------------------------------------------------------------
# 4. Collapse the multi-run cells down to single runs per the new
#    copy, and blank out the new 5th cell in the two data rows.
# ------------------------------------------------------------------
$d.Content.Find.Execute("1" + ".", $false, $false, $false, $false, $false, $true, 1, $false, "1.", 2) | Out-Null
$d.Content.Find.Execute("Test for displaying maze when list is empty. ", $false, $false, $false, $false, $false, $true, 1, $false, "Test for displaying maze when list is empty. ", 2) | Out-Null
$d.Content.Find.Execute("mazeList = [ ]", $false, $false, $false, $false, $false, $true, 1, $false, "mazeList = [ ]", 2) | Out-Null
$d.Content.Find.Execute("“No maze loaded.” Program will exit back to main menu. ", $false, $false, $false, $false, $false, $true, 1, $false, "“No maze loaded.” Program will exit back to main menu. ", 2) | Out-Null

$d.Content.Find.Execute("2" + ".", $false, $false, $false, $false, $false, $true, 1, $false, "2.", 2) | Out-Null
$d.Content.Find.Execute("Test for displaying maze when list is loaded. ", $false, $false, $false, $false, $false, $true, 1, $false, "Test for displaying maze when list is loaded. ", 2) | Out-Null
$d.Content.Find.Execute("mazeList = [“X”, “X”,”X”, “X”, “O” … …] ", $false, $false, $false, $false, $false, $true, 1, $false, "mazeList = [“X”, “X”,”X”, “X”, “O” … …] ", 2) | Out-Null

# ------------------------------------------------------------------
# 5. Remove the stray _GoBack bookmark from the closing paragraph; it
#    now lives in the new header cell instead.
# ------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()
